$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updated values for rows 2-7
$updates = @{
    2 = 2246
    3 = 641
    4 = 1651
    5 = 7635
    6 = 180
    7 = 211
}

# Both "展览" and "全部类型" sheets contain identical data and both need updating
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
